$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Commodity")
$ws.Activate()

# Insert a new blank row above row 6 so an empty line separates the
# "Dummy" commodity from the "Propulsion_of_Vehicles" block. This shifts
# all rows from 6 downward by one (old row 6 -> 7, ... old row 21 -> 22).
$ws.Rows.Item(6).Insert()

# The worksheet Table ("Tabelle3") does not auto-grow when rows are
# inserted via Rows.Insert(), so extend it (and its AutoFilter) to cover
# the new row, matching the updated data range A1:D22.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D22"))

$ws.Range("H12").Select()
